$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 198.4680276666666
$ws.Range("H2").Value = 595.4040829999999
$ws.Range("I2").Value = 0.2835009389723355
$ws.Range("J2").Value = 0.2835009389723355
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 131.4884046666667
$ws.Range("N2").Value = 394.465214
$ws.Range("O2").Value = 0.5054529031486338
$ws.Range("P2").Value = 0.5054529031486338
$ws.Range("Q2").Value = 26096.24433522986
$ws.Range("R2").Value = 234866.1990170687
$ws.Range("S2").Value = 0.1432963726489306
$ws.Range("T2").Value = 0.1432963726489306

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 198.4680276666666
$ws.Range("H3").Value = 595.4040829999999
$ws.Range("I3").Value = 0.2835009389723355
$ws.Range("J3").Value = 0.2835009389723355
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 37.83955
$ws.Range("N3").Value = 113.51865
$ws.Range("O3").Value = 0.1454585326350568
$ws.Range("P3").Value = 0.1454585326350568
$ws.Range("Q3").Value = 7509.940856294216
$ws.Range("R3").Value = 67589.46770664795
$ws.Range("S3").Value = 0.04123763058357671
$ws.Range("T3").Value = 0.04123763058357671

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 198.4680276666666
$ws.Range("H4").Value = 595.4040829999999
$ws.Range("I4").Value = 0.2835009389723355
$ws.Range("J4").Value = 0.2835009389723355
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 50.31467799999999
$ws.Range("N4").Value = 150.944034
$ws.Range("O4").Value = 0.1934140134300057
$ws.Range("P4").Value = 0.1934140134300057
$ws.Range("Q4").Value = 9985.854905343422
$ws.Range("R4").Value = 89872.6941480908
$ws.Range("S4").Value = 0.05483305441781453
$ws.Range("T4").Value = 0.05483305441781453

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 198.4680276666666
$ws.Range("H5").Value = 595.4040829999999
$ws.Range("I5").Value = 0.2835009389723355
$ws.Range("J5").Value = 0.2835009389723355
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 40.49714266666666
$ws.Range("N5").Value = 121.491428
$ws.Range("O5").Value = 0.1556745507863038
$ws.Range("P5").Value = 0.1556745507863038
$ws.Range("Q5").Value = 8037.388031188945
$ws.Range("R5").Value = 72336.4922807005
$ws.Range("S5").Value = 0.04413388132201367
$ws.Range("T5").Value = 0.04413388132201367

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 74.018453
$ws.Range("H6").Value = 222.055359
$ws.Range("I6").Value = 0.1057313924740739
$ws.Range("J6").Value = 0.1057313924740739
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 131.4884046666667
$ws.Range("N6").Value = 394.465214
$ws.Range("O6").Value = 0.5054529031486338
$ws.Range("P6").Value = 0.5054529031486338
$ws.Range("Q6").Value = 9732.568300864647
$ws.Range("R6").Value = 87593.11470778182
$ws.Range("S6").Value = 0.05344223927996825
$ws.Range("T6").Value = 0.05344223927996825

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 74.018453
$ws.Range("H7").Value = 222.055359
$ws.Range("I7").Value = 0.1057313924740739
$ws.Range("J7").Value = 0.1057313924740739
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 37.83955
$ws.Range("N7").Value = 113.51865
$ws.Range("O7").Value = 0.1454585326350568
$ws.Range("P7").Value = 0.1454585326350568
$ws.Range("Q7").Value = 2800.82495321615
$ws.Range("R7").Value = 25207.42457894535
$ws.Range("S7").Value = 0.01537953320274007
$ws.Range("T7").Value = 0.01537953320274007

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 74.018453
$ws.Range("H8").Value = 222.055359
$ws.Range("I8").Value = 0.1057313924740739
$ws.Range("J8").Value = 0.1057313924740739
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 50.31467799999999
$ws.Range("N8").Value = 150.944034
$ws.Range("O8").Value = 0.1934140134300057
$ws.Range("P8").Value = 0.1934140134300057
$ws.Range("Q8").Value = 3724.214628753133
$ws.Range("R8").Value = 33517.9316587782
$ws.Range("S8").Value = 0.02044993296395373
$ws.Range("T8").Value = 0.02044993296395373

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 74.018453
$ws.Range("H9").Value = 222.055359
$ws.Range("I9").Value = 0.1057313924740739
$ws.Range("J9").Value = 0.1057313924740739
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.49714266666666
$ws.Range("N9").Value = 121.491428
$ws.Range("O9").Value = 0.1556745507863038
$ws.Range("P9").Value = 0.1556745507863038
$ws.Range("Q9").Value = 2997.535851106961
$ws.Range("R9").Value = 26977.82265996265
$ws.Range("S9").Value = 0.01645968702741183
$ws.Range("T9").Value = 0.01645968702741183

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 405.8333793333334
$ws.Range("H10").Value = 1217.500138
$ws.Range("I10").Value = 0.5797112283523728
$ws.Range("J10").Value = 0.5797112283523728
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 131.4884046666667
$ws.Range("N10").Value = 394.465214
$ws.Range("O10").Value = 0.5054529031486338
$ws.Range("P10").Value = 0.5054529031486338
$ws.Range("Q10").Value = 53362.38360902217
$ws.Range("R10").Value = 480261.4524811996
$ws.Range("S10").Value = 0.2930167233585674
$ws.Range("T10").Value = 0.2930167233585674

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 405.8333793333334
$ws.Range("H11").Value = 1217.500138
$ws.Range("I11").Value = 0.5797112283523728
$ws.Range("J11").Value = 0.5797112283523728
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 37.83955
$ws.Range("N11").Value = 113.51865
$ws.Range("O11").Value = 0.1454585326350568
$ws.Range("P11").Value = 0.1454585326350568
$ws.Range("Q11").Value = 15356.55244895264
$ws.Range("R11").Value = 138208.9720405737
$ws.Range("S11").Value = 0.08432394462820247
$ws.Range("T11").Value = 0.08432394462820247

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 405.8333793333334
$ws.Range("H12").Value = 1217.500138
$ws.Range("I12").Value = 0.5797112283523728
$ws.Range("J12").Value = 0.5797112283523728
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 50.31467799999999
$ws.Range("N12").Value = 150.944034
$ws.Range("O12").Value = 0.1934140134300057
$ws.Range("P12").Value = 0.1934140134300057
$ws.Range("Q12").Value = 20419.37580280852
$ws.Range("R12").Value = 183774.3822252767
$ws.Range("S12").Value = 0.1121242753060709
$ws.Range("T12").Value = 0.1121242753060709

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 405.8333793333334
$ws.Range("H13").Value = 1217.500138
$ws.Range("I13").Value = 0.5797112283523728
$ws.Range("J13").Value = 0.5797112283523728
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 40.49714266666666
$ws.Range("N13").Value = 121.491428
$ws.Range("O13").Value = 0.1556745507863038
$ws.Range("P13").Value = 0.1556745507863038
$ws.Range("Q13").Value = 16435.09226175745
$ws.Range("R13").Value = 147915.830355817
$ws.Range("S13").Value = 0.09024628505953203
$ws.Range("T13").Value = 0.09024628505953203

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 21.741411
$ws.Range("H14").Value = 65.224233
$ws.Range("I14").Value = 0.03105644020121776
$ws.Range("J14").Value = 0.03105644020121776
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 131.4884046666667
$ws.Range("N14").Value = 394.465214
$ws.Range("O14").Value = 0.5054529031486338
$ws.Range("P14").Value = 0.5054529031486338
$ws.Range("Q14").Value = 2858.743447592318
$ws.Range("R14").Value = 25728.69102833086
$ws.Range("S14").Value = 0.01569756786116746
$ws.Range("T14").Value = 0.01569756786116746

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 21.741411
$ws.Range("H15").Value = 65.224233
$ws.Range("I15").Value = 0.03105644020121776
$ws.Range("J15").Value = 0.03105644020121776
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 37.83955
$ws.Range("N15").Value = 113.51865
$ws.Range("O15").Value = 0.1454585326350568
$ws.Range("P15").Value = 0.1454585326350568
$ws.Range("Q15").Value = 822.68520860505
$ws.Range("R15").Value = 7404.16687744545
$ws.Range("S15").Value = 0.004517424220537522
$ws.Range("T15").Value = 0.004517424220537522

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 21.741411
$ws.Range("H16").Value = 65.224233
$ws.Range("I16").Value = 0.03105644020121776
$ws.Range("J16").Value = 0.03105644020121776
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 50.31467799999999
$ws.Range("N16").Value = 150.944034
$ws.Range("O16").Value = 0.1934140134300057
$ws.Range("P16").Value = 0.1934140134300057
$ws.Range("Q16").Value = 1093.912093730658
$ws.Range("R16").Value = 9845.208843575921
$ws.Range("S16").Value = 0.0060067507421665
$ws.Range("T16").Value = 0.0060067507421665

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 21.741411
$ws.Range("H17").Value = 65.224233
$ws.Range("I17").Value = 0.03105644020121776
$ws.Range("J17").Value = 0.03105644020121776
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 40.49714266666666
$ws.Range("N17").Value = 121.491428
$ws.Range("O17").Value = 0.1556745507863038
$ws.Range("P17").Value = 0.1556745507863038
$ws.Range("Q17").Value = 880.4650230416358
$ws.Range("R17").Value = 7924.185207374723
$ws.Range("S17").Value = 0.004834697377346281
$ws.Range("T17").Value = 0.004834697377346281
